$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value of 45182 (2023-09-13)
# for every data row (2..385). Update it to 45184 (2023-09-15).
for ($r = 2; $r -le 385; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}
